# Commit: Add write K_polinoms in ROM
# Update K_polinoms lookup table values (rows 2-9, Sheet1) to the refreshed
# polynomial coefficients.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.31
$ws.Range("C2").Value = 1.28
$ws.Range("D2").Value = 1.29
$ws.Range("E2").Value = 1.3
$ws.Range("I2").Value = 1.32
$ws.Range("J2").Value = 1.29
$ws.Range("K2").Value = 1.36
$ws.Range("N2").Value = 1.28
$ws.Range("O2").Value = 1.32
$ws.Range("P2").Value = 1.4
$ws.Range("Q2").Value = 1.34
$ws.Range("R2").Value = 1.29
$ws.Range("S2").Value = 1.31
$ws.Range("U2").Value = 1.42
$ws.Range("W2").Value = 1.26
$ws.Range("X2").Value = 1.33
$ws.Range("AA2").Value = 1.28
$ws.Range("AB2").Value = 1.3
$ws.Range("AD2").Value = 1.32
$ws.Range("AE2").Value = 1.28
$ws.Range("AF2").Value = 1.31
$ws.Range("AG2").Value = 1.3

# Row 3
$ws.Range("F3").Value = 1.14
$ws.Range("G3").Value = 1.11
$ws.Range("H3").Value = 1.12
$ws.Range("J3").Value = 1.13
$ws.Range("V3").Value = 1.13
$ws.Range("Y3").Value = 1.16
$ws.Range("AD3").Value = 1.14

# Row 4
$ws.Range("E4").Value = 1.09
$ws.Range("N4").Value = 1.07
$ws.Range("W4").Value = 1.09
$ws.Range("X4").Value = 1.1
$ws.Range("AD4").Value = 1.1
$ws.Range("AF4").Value = 1.09

# Row 5
$ws.Range("D5").Value = 1.07
$ws.Range("U5").Value = 1.1
$ws.Range("V5").Value = 1.07
$ws.Range("X5").Value = 1.08
$ws.Range("AA5").Value = 1.07
$ws.Range("AD5").Value = 1.08
$ws.Range("AG5").Value = 1.07

# Row 6
$ws.Range("I6").Value = 1.06
$ws.Range("M6").Value = 1.06
$ws.Range("P6").Value = 1.07
$ws.Range("T6").Value = 1.07
$ws.Range("U6").Value = 1.07
$ws.Range("Z6").Value = 1.06
$ws.Range("AB6").Value = 1.07

# Row 7
$ws.Range("B7").Value = 1.05
$ws.Range("I7").Value = 1.05
$ws.Range("J7").Value = 1.05
$ws.Range("N7").Value = 1.03
$ws.Range("AE7").Value = 1.06

# Row 9
$ws.Range("O9").Value = 1.05
$ws.Range("T9").Value = 1.06
